$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.956.56'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.194.23'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.86%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.71'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.80'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +2.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.34'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.113'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.14%  '
$ws.Range("E11").Value = '  +2.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.745.43'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.90%  '
$ws.Range("E13").Value = '  -1.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.11'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000174'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.996.22'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.184.76'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.23'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.11'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.40'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '382.42'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("E23").Value = '  +3.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.25'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("E25").Value = '  +2.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.87'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +12.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0900'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("E29").Value = '  +1.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.37'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.51%  '
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.37'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.63%  '
$ws.Range("E33").Value = '  +3.50%  '
$ws.Range("E34").Value = '  +4.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.27'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.761.37'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +5.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.57'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0712'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.77%  '
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.74'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("E43").Value = '  +4.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0286'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +5.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.235.05'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.19%  '
$ws.Range("E47").Value = '  +1.93%  '
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.798'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +6.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.50'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("E51").Value = '  -0.04%  '
